$d = $word.ActiveDocument

# 1. Title
$d.Content.Find.Execute('Review 195: Can Mamba Learn How to Learn? A Comparative Study on In-Context Learning Tasks', $true, $false, $false, $false, $false, $true, 1, $false, 'Review 194: Mamba: Linear-Time Sequence Modeling with Selective State Spaces', 2)

# 2. Paper link (bold)
$d.Content.Find.Execute('Paper: https://arxiv.org/abs/2402.04248v2', $true, $false, $false, $false, $false, $true, 1, $false, 'Paper: https://arxiv.org/abs/2507.06204v1', 2)

# 3. arxiv abs link
$d.Content.Find.Execute('https://arxiv.org/abs/2402.04248', $true, $false, $false, $false, $false, $true, 1, $false, 'https://arxiv.org/abs/2312.00752', 2)

# 4. Paragraph 1 (intro)
$d.Content.Find.Execute('אוקיי, סוקרים מאמר הבא בסדרת ממבה (מה שבא אחרי). בניגוד להצהרותיי בסוף סקירתי הקודמת לא תהיי Mamba MoE אלא מאמר אחר. הסיבה היא שלדעתי כמות המאמרים על MoE היא גדולה מאוד והמאמר הזה רק מציע להלביש אותו על Mamba ללא חידושים מעניינים אחרים אז החלטתי לדלג.', $true, $false, $false, $false, $false, $true, 1, $false, 'זה קורה עכשיו, אחרי 9 סקירות שחלקם היו די לא פשוטות הגענו למטרתנו הקדושה שזה Mamba. מכיוון שאני מפרסם סקירות בשלשות (באתר MDLI) אני אוציא עוד 2 סקירות נוספות של שכלולי ממבה (אחד מהם Mamba MoE והשני עוד הוחלט).', 2)

# 5. Paragraph 2
$d.Content.Find.Execute('המאמר שנסקור היום בודק את האם מודלים המבוססים על ארכיטקטורת ממבה על למידת in-context (או ICL). למעשה ICL היא יכולת של מודל לבצע למידה על בסיס כמה דוגמאות בודדות (גם נקרא למידת few-shot) ללא שינוי של משקלי המודל. בגדול יכולת זו של הטרנספורמרים לא מאוד מפתיעה  כי ה״חיזויים״ שלהם תלויים ביחסים בין חלקי הדאטה השונים(טוקנים) באופן מפורש באמצעות מנגנון ה-attention שלהם. כמובן יש מחקרים לא מעטים ומעניינים שחוקרים את התופעה המרתקת הזו ואני ממליץ לכם בחום להעיף מבט.', $true, $false, $false, $false, $false, $true, 1, $false, 'האמת שאחרי שאנו הבנו מה- (SSM (space-state models ואיך ניתן לבנות ארכיטקטורה מבוססת עליהם לעיבוד דאטה סדרתי, השכלול המוצע על ידי mamba הוא די אינטואיטיבי ומתבקש. כמו שאתם זוכרים SSM ממומשת בתור מערכת דינמית(DLS) לינארית כאשר הקלט למערכת זו היא ייצוג וקטורי (embeddings) של איברי הסדרה (= טוקנים). ', 2)

# 6. Paragraph 3
$d.Content.Find.Execute('לעומת זאת הארכיטקטורה של ממבה לא לוקחת את היחסים בין הטוקנים השונים של הדאטה באופן מפורש ודוחסת את ה״עבר״ בוקטור אחד אז היכולת שלה לבצע ICL היא פחות אינטואיטיבי.  זה אכן פחות קורה. המאמר בדק כמה ארכיטקטורות מבוססות (SSM (state-space models כמו ממבה, S4 וגם S4-mamba ללא מנגנון attention של הטרנספורמרים והשוו את יכולות ICL שלהם עם ארכיטקטורות הברידיות: כלומר שילוב של ממבה עם מנגנון ה-attention של הטרנספורמרים.', $true, $false, $false, $false, $false, $true, 1, $false, 'בשלב הראשון המערכת הדינמית מחשבת וקטור s הוא הוא ייצוג דחוס של זיכרון כלומר וקטור ״הזוכר״ את המידע הרלוונטי עבור כל הטוקנים הקודמים לטוקן הנוכחי. בשלב השני מחשבים את הפלט עבור טוקן זה המוזן לשכבה הבאה (שיכולה להיות גם שכבת שמייצרת פלט סופי). כל חישובים אלו מתבצעים באמצעות מיפויים לינאריים כלומר מכפלות במטריצות. חשוב להבין שכל המעברים בין ייצוגי הזכרון בין הטוקנים הם לינאריים ונשלטים על ידי אותה מטריצה A ווקטורי B, C וסקלר delta. ד״א פרמטר delta מגדיר (באופן גס) את קצב דעיכה של הזכרון (כלומר ככל ש delta גבוה יותר אנו נוטים ״לזכור״ פחות מהטוקנים הקודמים).', 2)

# 7. Paragraph 4 + two new paragraphs after it
$d.Content.Find.Execute('איך משלבים ממבה עם הטרנספורמר? המאמר בדק שתי גישות (די דומות). בגישה הראשונה הוא החליף את MLP שיש בבלוקי טרנספורמר אחרי ה-attention במנגנון של ממבה. הגישה השנייה (הכי מוצלחת) הנקראת MambaFormer בנוסף מחליפה את הקידוד המיקומי (positional encoding) של עוברים הטוקנים בממבה נוספת. ', $true, $false, $false, $false, $false, $true, 1, $false, 'מה היתרונות של הארכיטקטורה הזו? היא בעלת תכונה הדואליות המיוחלת המשלבת 2 התכונות הבאות:', 2)
$p4 = $d.Paragraphs.Item(12)
$p4.Range.InsertParagraphAfter()
$p4a = $d.Paragraphs.Item(13)
$p4a.Range.Text = 'ניתן לחזות באופן מקבילי (בו זמנית) כמה טוקנים במהלך אימון (כמו בטרנספורמרים)'
$p4a.Range.InsertParagraphAfter()
$p4b = $d.Paragraphs.Item(14)
$p4b.Range.Text = 'חיזוי מהיר של טוקן במהלך היסק (ללא התחשבות בכל הטוקנים בחלון ההקשר כמו בטרנספורמרים שמביא לנו את הסיבוכיות הריבועית).'

# 8. Paragraph 5 (was index 14, now 16 after 2 insertions)
$d.Content.Find.Execute('כאמור MambaFormer הגיע לביצועים הטובים ביותר מכל הארכיטקטורות הלא היברידיות (הטרנספורמר הטהור וכמה וריאנטים של SSM) באופן לא מפתיע בכלל. הרי MLP (רק 2 שכבות) ממדלים הפעולה די פשוטה ו-mamba היא למעשה מנגנון של זכרון הדוחס את המידע המהותי (בתקווה) של העבר (בטוקנים הקודמים). לא פלא שזה ניצח את כולם.', $true, $false, $false, $false, $false, $true, 1, $false, 'כלומר הארכיטקטורות מסוג זה הם יעילות בזמן האימון ומהירות בזמן ההיסק. אבל כמו שאתם יכולים לנחש יש לנו מחיר לשלם על כל התכונות הנחמדות האלו. ומחיר הוא כמובן יכולת של  המודל למדל תלויות מורכבות של הדאטה. עקב כך מאמרים כמו Hyena, H3, S4 ניסו ניסו לבנות את הפרמטרים של DLS (המגדירה מעברים בין ייצוגי הזכרון ויצירת הפלט) בצורה חכמה (ודי מורכבת).', 2)

# 9. Paragraph 6
$d.Content.Find.Execute('נשאר לנו רק לציין איזה משימות ניתנו למודלים אלו כדי לבחון את יכולות ICL שלהם. אחת המשימות היא לתת למודל כמה זוגות של (x, f(x)) עבור פונקציה f לינארית ולבקש ממנו לחשב (f(x עבור x-ים נוספים. משימה אחרת היתה לתת לה נקודות שנדגמו מ Gaussian Mixture מסוים ולבקש ממנו לדגום עוד נקודות. טבלה עם כל המשימות מצורפת לפוסט.', $true, $false, $false, $false, $false, $true, 1, $false, 'אבל מתברר שזה לא מספיק. מעברים לינאריים עם פרמטרי DLS קבועים לא מסוגלת למדל דאטה מורכב (כמו שפה טבעית). אחד המשימות שמודל כזה נכשל עליה הוא העתקת טוקנים הבאים אחרי טוקן ספציפי (וזה די הגיוני לאור הפרמטרים הקבועים של DLS). כמו שאתם יכולים כבר לנחש אולי מחברי ממבה מציעים לעשות חלק מהפרמטרים (B, C ו delta) תלוים בייצוג הטוקן הנוכחי. התלות הזו היא לינארית עם מטריצות נלמדות. וזה עוזר לנו להתחשב בפיסת הקלט הנוכחית בצורה יותר טובה. כאמור B, C מגדירות את האופן בו ייצוג הזכרון והפלט עבור הטוקן הנוכחי בהתאמה ואז יש לנו סיכוי יותר טוב להצליח במשימות מהסוג שתיארתי לפני. בנוסף תלות של delta בייצוג הטוקן הנוכחי מקנה לנו אפשרות לשחק עם קצב דעיכה בצורה יותר גרנולרית שמקנה לנו יכולת ״לשכוח״ ו״לזכור״ איפה שצריך. ', 2)

# 10. Paragraph 7 + insert empty/text/empty/text after it
$d.Content.Find.Execute('נתראה בסקירה ממבה הבאה והאחרונה (לא בחרתי עדיין).', $true, $false, $false, $false, $false, $true, 1, $false, 'אבל האם איבדנו את הדואליות שלנו בדרך. מתברר שלא, הרי המעבר בין ייצוגי הזיכרון של הטוקנים עדיין לא תלוי במיקום של הטוקן אלא בייצוג. כלומר אנו עדיין יכולים לחזות מספר טוקנים בו זמנית כי אנו יכולים לחשב את כל הפרמטרים מראש (לא צריך לחשב את המצב הזכרון הקודם באופן מפורש). וכמובן אין צורך להתחשב בכל הטוקנים בתוך חלון ההקשר במהלך ההיסק כי הזיכרון עדיין מיוצג על יד וקטור אחר. אז יש דואליות!', 2)
$p7 = $d.Paragraphs.Item(20)
$p7.Range.InsertParagraphAfter()
$p7empty1 = $d.Paragraphs.Item(21)
$p7empty1.Range.InsertParagraphAfter()
$p7a = $d.Paragraphs.Item(22)
$p7a.Range.Text = 'מה שכן קורה הוא זה החישובים הופכים לקצת יותר מורכבים (שימו לב שהחלק הכי בעייתי בחישוב שהוא העלאה של מטריצה בחזקה לא השתנה כי A נותרה קבועה). בסוף המאמר מציע כמה שכלולים לאופן חישוב המייעלים ומזרזים אותו (המשחקים בין זכרון מהיר ואיטי של GPU).'
$p7a.Range.InsertParagraphAfter()
$p7empty2 = $d.Paragraphs.Item(23)
$p7empty2.Range.InsertParagraphAfter()
$p7b = $d.Paragraphs.Item(24)
$p7b.Range.Text = 'זה וזה עכשיו אתם יודעים מה זה ממבה. נתראה ב Mamba MoE עוד כמה ימים.'

Write-Output "Edit complete"
